# Updates cryptos list values (Price / Volume(1h) columns) per the
# Mon Oct  9 18:38:19 UTC 2023 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.716.83'
$ws.Range("E2").Value = '  -0.54%  '
$ws.Range("D3").Value = '1.584.09'
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = "'207.43"
$ws.Range("E5").Value = '  -1.97%  '
$ws.Range("E6").Value = '  -3.12%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = "'22.36"
$ws.Range("E8").Value = '  -4.39%  '
$ws.Range("E9").Value = '  -1.11%  '
$ws.Range("E10").Value = '  -3.40%  '
$ws.Range("E11").Value = '  -1.69%  '
$ws.Range("D12").Value = '1.809.82'
$ws.Range("E12").Value = '  -2.96%  '
$ws.Range("D13").Value = '1.564.20'
$ws.Range("E13").Value = '  -4.35%  '
$ws.Range("D14").Value = "'3.87"
$ws.Range("E14").Value = '  -3.88%  '
$ws.Range("D15").Value = "'0.532"
$ws.Range("E15").Value = '  -5.13%  '
$ws.Range("D16").Value = '27.686.60'
$ws.Range("E16").Value = '  -0.69%  '
$ws.Range("D17").Value = "'63.02"
$ws.Range("E17").Value = '  -3.45%  '
$ws.Range("D18").Value = "'218.01"
$ws.Range("E18").Value = '  -4.74%  '
$ws.Range("E19").Value = '  -4.33%  '
$ws.Range("D20").Value = '0.0₃0693'
$ws.Range("E20").Value = '  -3.55%  '
$ws.Range("E21").Value = '  +0.18%  '
$ws.Range("E22").Value = '  -4.13%  '
$ws.Range("D23").Value = "'9.57"
$ws.Range("E23").Value = '  -4.59%  '
$ws.Range("D24").Value = "'1.97"
$ws.Range("E24").Value = '  -5.03%  '
$ws.Range("D25").Value = "'153.62"
$ws.Range("E25").Value = '  -1.02%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D27").Value = "'6.71"
$ws.Range("E27").Value = '  -2.37%  '
$ws.Range("D28").Value = "'15.09"
$ws.Range("E28").Value = '  -2.93%  '
$ws.Range("E29").Value = '  -4.26%  '
$ws.Range("E30").Value = '  -1.92%  '
$ws.Range("E31").Value = '  -3.43%  '
$ws.Range("E32").Value = '  -4.87%  '
$ws.Range("D33").Value = '1.377.40'
$ws.Range("E33").Value = '  -1.17%  '
$ws.Range("E34").Value = '  -4.99%  '
$ws.Range("E35").Value = '  -4.93%  '
$ws.Range("D36").Value = "'0.970"
$ws.Range("E36").Value = '  -4.73%  '
$ws.Range("D37").Value = "'2.31"
$ws.Range("E37").Value = '  -1.67%  '
$ws.Range("E38").Value = '  -3.27%  '
$ws.Range("D39").Value = "'0.538"
$ws.Range("E39").Value = '  -3.72%  '
$ws.Range("D40").Value = "'0.816"
$ws.Range("E40").Value = '  -3.74%  '
$ws.Range("D42").Value = "'0.977"
$ws.Range("E42").Value = '  -3.43%  '
$ws.Range("E43").Value = '  -1.73%  '
$ws.Range("E44").Value = '  +1.81%  '
$ws.Range("D45").Value = "'63.68"
$ws.Range("E45").Value = '  -3.28%  '
$ws.Range("D46").Value = "'5.21"
$ws.Range("E46").Value = '  -4.14%  '
$ws.Range("D47").Value = '1.719.96'
$ws.Range("E47").Value = '  -3.03%  '
$ws.Range("D48").Value = "'87.80"
$ws.Range("E48").Value = '  -0.99%  '
$ws.Range("D49").Value = '0.0₇0998'
$ws.Range("E49").Value = '  -3.03%  '
$ws.Range("D50").Value = "'0.0975"
$ws.Range("E50").Value = '  -4.65%  '
$ws.Range("E51").Value = '  -1.73%  '
